# Upload new version with timestamp
# Fills the previously-empty template row (row 4) with 14 product rows,
# inserts the extra rows needed, recomputes the totals row and pushes the
# footer row down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert 13 new rows right after the existing template row
#    (row 4). This pushes the old "totals" row (5) down to row 18 and the
#    old "footer" row (6) down to row 19, carrying their formatting with
#    them. The newly inserted rows 5-17 inherit row 4's formatting.
# ---------------------------------------------------------------------
$ws.Range("A5:A17").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. The product data to fill into rows 4-17.
#    Columns: A = index, B = product name, H = balance, L = price, N = count
# ---------------------------------------------------------------------
$data = @(
    @(1,  "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS.",       "0:0",    114,   1),
    @(2,  "AUGMENTIN 457MG/5ML SUSP. 70 ML",           "1:0",    137,   1),
    @(3,  "BLOKATENS 10/160MG 28 F.C.TABS.",           "0:0",    160,   1),
    @(4,  "COLOVATIL 30 F.C. TABS",                    "0:0",    63,    1),
    @(5,  "GAVISCON LIQUID 24 SACHETS 10 ML",          "0:9",    12,    0.04),
    @(6,  "GINKGO BILOBA 30 CAPS.",                    "0:0",    186,   1),
    @(7,  "MILGA ADVANCE 30 F.C. TABS",                "0:0",    136.5, 1),
    @(8,  "PERLOC 40MG 14 F.C.TAB.",                   "0:0",    68.25, 1),
    @(9,  "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML",  "2:0",    18,    1),
    @(10, "RIVO 320MG 20*10 TABS",                     "1:2",    14.1,  0.1),
    @(11, "VASTAREL MR 35MG 30 F.C.TAB.",               "2:0",    175,   1),
    @(12, "WATER FOR INJECTION AMP. 5 ML",              "7816:0", 2.5,   1),
    @(13, "سويت كوكو",                                  "22:0",   25,    1),
    @(14, "مرطب شفاه لونا جوز هند ابيض",                  "3:0",    20,    1)
)

$rowHeights = @{
    4=24.75; 5=25.5; 6=24.75; 7=25.5; 8=25.5; 9=24.75; 10=25.5; 11=24.75;
    12=25.5; 13=25.5; 14=24.75; 15=25.5; 16=24.75; 17=25.5
}

# The B:G and H:K merged blocks must hold text (product name / balance
# ratio) so force a text number format before writing the strings -
# otherwise values such as "0:0" or "1:2" would be reinterpreted as times.
$ws.Range("B4:G17").NumberFormat = "@"
$ws.Range("H4:K17").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 4 + $i
    $item = $data[$i]

    $ws.Cells.Item($r, 1).Value2 = $item[0]        # A - index
    $ws.Cells.Item($r, 2).Value2 = $item[1]        # B - product name
    $ws.Cells.Item($r, 8).Value2 = $item[2]        # H - balance
    $ws.Cells.Item($r, 12).Value2 = $item[3]       # L - price
    $ws.Cells.Item($r, 14).Value2 = $item[4]       # N - transactions count

    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# ---------------------------------------------------------------------
# 3. Totals row (now row 18): sum of the price column.
# ---------------------------------------------------------------------
$ws.Range("K18").Value2 = 1131.3499999999999

Write-Host "Workbook updated."
